# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G (header "K") previously held the old "Strike#" style value; it is
# regenerated here to hold the newly-computed K (strike count) values per
# row, row 2 through row 25 of Sheet1 (row 24 is unchanged by the regen).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 4
    4  = 5
    5  = 4
    6  = 6
    7  = 5
    8  = 7
    9  = 4
    10 = 4
    11 = 3
    12 = 3
    13 = 5
    14 = 2
    15 = 5
    16 = 3
    17 = 3
    18 = 5
    19 = 4
    20 = 4
    21 = 3
    22 = 7
    23 = 2
    25 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
